# The codeforiati:category-code and codeforiati:category-name columns
# (F and G) were swapped: what used to be in column F now belongs in
# column G and vice-versa, for every row of the sheet (including the
# header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 6).End(-4162).Row  # xlUp = -4162

# The values that move into column G are short numeric-looking codes
# (e.g. "111"). Format that range as Text first so Excel keeps them as
# text instead of auto-converting them to numbers, matching the original
# column's text values.
$ws.Range("G2:G" + $lastRow).NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2
    $fCell.Value2 = $gVal
    $gCell.Value2 = $fVal
}
